# Add the new "Gaussian Olympic" sheet at the end of the workbook (after the
# last existing sheet), matching sheetId=3 / rId3 ordering in the target file.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Gaussian Olympic"

# Column A is widened to fit the long event labels.
$ws.Columns.Item(1).ColumnWidth = 40.140625

# Seed the shared-string table in the same order new strings first appear in
# the authored workbook: the two section headers and the training note are
# entered before the data rows are filled in.
$ws.Range("A3").Value = "MSEs:"
$ws.Range("A19").Value = "Average logp on hold out:"
$ws.Range("A1").Value = "Training was first 20 data points"

# --- MSEs section -------------------------------------------------------
$ws.Range("A4").Value = "100 metres, Men"
$ws.Range("B4").Value = 0.0067224457779073996

$ws.Range("A5").Value = "400 metres, Men"
$ws.Range("B5").Value = 0.16619999999999599

$ws.Range("A6").Value = "800 metres, Men"
$ws.Range("B6").Value = 1.10737252809609

$ws.Range("A7").Value = "110 metres Hurdles, Men : 0.06377592219848145"
$ws.Range("B7").Value = 0.063775922198481397

$ws.Range("A8").Value = "200 metres, Men"
$ws.Range("B8").Value = 0.16170381418090601

$ws.Range("A9").Value = "400 metres Hurdles, Men"
$ws.Range("B9").Value = 0

$ws.Range("A10").Value = "1,500 metres, Men"
$ws.Range("B10").Value = 0.176333333333339

$ws.Range("A11").Value = "5,000 metres, Men"
$ws.Range("B11").Value = 45.279674999999699

$ws.Range("A12").Value = "4 x 100 metres Relay, Men"
$ws.Range("B12").Value = 0.156799999999997

$ws.Range("A13").Value = "4 x 400 metres Relay, Men"
$ws.Range("B13").Value = 0.091875000000011905

$ws.Range("A14").Value = "3,000 metres Steeplechase, Men"
$ws.Range("B14").Value = 1.6576333333333599

$ws.Range("A15").Value = "100 metres, Women"
$ws.Range("B15").Value = 0

# --- Average logp on hold out section ----------------------------------
$ws.Range("A20").Value = "100 metres, Men"
$ws.Range("B20").Value = -7.1487083884933096

$ws.Range("A21").Value = "400 metres, Men"
$ws.Range("B21").Value = -10.7229850624663

$ws.Range("A22").Value = "800 metres, Men"
$ws.Range("B22").Value = -12.0064475038676

$ws.Range("A23").Value = "110 metres Hurdles, Men"
$ws.Range("B23").Value = -9.7046885047827995

$ws.Range("A24").Value = "200 metres, Men"
$ws.Range("B24").Value = -6.7996248520653699

$ws.Range("A25").Value = "400 metres Hurdles, Men"
$ws.Range("B25").Value = -3.1620820157848599

$ws.Range("A26").Value = "1,500 metres, Men"
$ws.Range("B26").Value = -10.1513154439119

$ws.Range("A27").Value = "5,000 metres, Men"
$ws.Range("B27").Value = -15.953386773586899

$ws.Range("A28").Value = "4 x 100 metres Relay, Men"
$ws.Range("B28").Value = -5.5655457085016602

$ws.Range("A29").Value = "4 x 400 metres Relay, Men"
$ws.Range("B29").Value = -8.6598274347477897

$ws.Range("A30").Value = "3,000 metres Steeplechase, Men"
$ws.Range("B30").Value = -13.426897982511599

$ws.Range("A31").Value = "100 metres, Women"
$ws.Range("B31").Value = -1.7103210301987599

# New sheet becomes the active/selected tab, with row 4 selected (matches the
# authored sheetView's activeCell/sqref).
[void]$ws.Rows.Item(4).Select()
